$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il11"
$ws.Range("C2").Value = "Il11ra1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04156866666666666
$ws.Range("H2").Value = 0.124706
$ws.Range("I2").Value = 0.02383845226880681
$ws.Range("J2").Value = 0.02383845226880681
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.440362666666666
$ws.Range("N2").Value = 13.321088
$ws.Range("O2").Value = 0.08107461047911722
$ws.Range("P2").Value = 0.08107461047911722
$ws.Range("Q2").Value = 0.1845799555697777
$ws.Range("R2").Value = 1.661219600128
$ws.Range("S2").Value = 0.00193269323211854
$ws.Range("T2").Value = 0.00193269323211854

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il11"
$ws.Range("C3").Value = "Il11ra1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04156866666666666
$ws.Range("H3").Value = 0.124706
$ws.Range("I3").Value = 0.02383845226880681
$ws.Range("J3").Value = 0.02383845226880681
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 46.01708333333334
$ws.Range("N3").Value = 138.05125
$ws.Range("O3").Value = 0.8402054937183234
$ws.Range("P3").Value = 0.8402054937183234
$ws.Range("Q3").Value = 1.912868798055556
$ws.Range("R3").Value = 17.2158191825
$ws.Range("S3").Value = 0.02002919855799351
$ws.Range("T3").Value = 0.02002919855799351

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il11"
$ws.Range("C4").Value = "Il11ra1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04156866666666666
$ws.Range("H4").Value = 0.124706
$ws.Range("I4").Value = 0.02383845226880681
$ws.Range("J4").Value = 0.02383845226880681
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.311397666666667
$ws.Range("N4").Value = 12.934193
$ws.Range("O4").Value = 0.07871989580255943
$ws.Range("P4").Value = 0.07871989580255942
$ws.Range("Q4").Value = 0.1792190524731111
$ws.Range("R4").Value = 1.612971472258
$ws.Range("S4").Value = 0.001876560478694759
$ws.Range("T4").Value = 0.001876560478694758

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Il11"
$ws.Range("C5").Value = "Il11ra1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.218041
$ws.Range("H5").Value = 3.654123
$ws.Range("I5").Value = 0.6985119939686074
$ws.Range("J5").Value = 0.6985119939686074
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.440362666666666
$ws.Range("N5").Value = 13.321088
$ws.Range("O5").Value = 0.08107461047911722
$ws.Range("P5").Value = 0.08107461047911722
$ws.Range("Q5").Value = 5.408543782869334
$ws.Range("R5").Value = 48.676894045824
$ws.Range("S5").Value = 0.05663158782599632
$ws.Range("T5").Value = 0.05663158782599632

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Il11"
$ws.Range("C6").Value = "Il11ra1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.218041
$ws.Range("H6").Value = 3.654123
$ws.Range("I6").Value = 0.6985119939686074
$ws.Range("J6").Value = 0.6985119939686074
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 46.01708333333334
$ws.Range("N6").Value = 138.05125
$ws.Range("O6").Value = 0.8402054937183234
$ws.Range("P6").Value = 0.8402054937183234
$ws.Range("Q6").Value = 56.05069420041668
$ws.Range("R6").Value = 504.45624780375
$ws.Range("S6").Value = 0.5868936147605643
$ws.Range("T6").Value = 0.5868936147605643

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Il11"
$ws.Range("C7").Value = "Il11ra1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.218041
$ws.Range("H7").Value = 3.654123
$ws.Range("I7").Value = 0.6985119939686074
$ws.Range("J7").Value = 0.6985119939686074
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.311397666666667
$ws.Range("N7").Value = 12.934193
$ws.Range("O7").Value = 0.07871989580255943
$ws.Range("P7").Value = 0.07871989580255942
$ws.Range("Q7").Value = 5.251459125304335
$ws.Range("R7").Value = 47.26313212773901
$ws.Range("S7").Value = 0.0549867913820468
$ws.Range("T7").Value = 0.05498679138204679

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Il11"
$ws.Range("C8").Value = "Il11ra1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.4841556666666667
$ws.Range("H8").Value = 1.452467
$ws.Range("I8").Value = 0.2776495537625858
$ws.Range("J8").Value = 0.2776495537625858
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.440362666666666
$ws.Range("N8").Value = 13.321088
$ws.Range("O8").Value = 0.08107461047911722
$ws.Range("P8").Value = 0.08107461047911722
$ws.Range("Q8").Value = 2.149826747121777
$ws.Range("R8").Value = 19.348440724096
$ws.Range("S8").Value = 0.02251032942100236
$ws.Range("T8").Value = 0.02251032942100236

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Il11"
$ws.Range("C9").Value = "Il11ra1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.4841556666666667
$ws.Range("H9").Value = 1.452467
$ws.Range("I9").Value = 0.2776495537625858
$ws.Range("J9").Value = 0.2776495537625858
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 46.01708333333334
$ws.Range("N9").Value = 138.05125
$ws.Range("O9").Value = 0.8402054937183234
$ws.Range("P9").Value = 0.8402054937183234
$ws.Range("Q9").Value = 22.27943165930556
$ws.Range("R9").Value = 200.51488493375
$ws.Range("S9").Value = 0.2332826803997655
$ws.Range("T9").Value = 0.2332826803997655

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Il11"
$ws.Range("C10").Value = "Il11ra1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4841556666666667
$ws.Range("H10").Value = 1.452467
$ws.Range("I10").Value = 0.2776495537625858
$ws.Range("J10").Value = 0.2776495537625858
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.311397666666667
$ws.Range("N10").Value = 12.934193
$ws.Range("O10").Value = 0.07871989580255943
$ws.Range("P10").Value = 0.07871989580255942
$ws.Range("Q10").Value = 2.087387611570111
$ws.Range("R10").Value = 18.786488504131
$ws.Range("S10").Value = 0.02185654394181788
$ws.Range("T10").Value = 0.02185654394181787
